# Add upstream/downstream adapter columns to the sequences worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("C1").Value = "upstream adapter"
$ws.Range("D1").Value = "downstream adapter"

# MinE row
$ws.Range("C2").Value = "caaacgcttgttcggaggataagtt"
$ws.Range("D2").Value = "gcccgctgtaaaagcgca"

# MinD row
$ws.Range("C3").Value = "tgatccctttttaacaaggaatttct"
$ws.Range("D3").Value = "gttatggcattactcgatttctttc"

# Column widths (values chosen so the stored OOXML width matches the
# target column widths of 38.1640625 / 25.58203125 / 29.58203125 once the
# runtime quantizes ColumnWidth to whole-pixel steps).
$ws.Columns.Item(2).ColumnWidth = 37.4297
$ws.Columns.Item(3).ColumnWidth = 24.8582
$ws.Columns.Item(4).ColumnWidth = 28.8582

# Match the saved selection/active cell
$ws.Range("D3").Select()
